$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (values are stored as text in the sheet).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated coin data (Price and Volume columns, plus the
# Name/Link columns for rows whose ranking shifted).
$ws.Range("D2").Value = '31.163.96'
$ws.Range("E2").Value = '  +1.88%  '
$ws.Range("D3").Value = '1.989.31'
$ws.Range("E3").Value = '  +5.65%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '0.7903'
$ws.Range("E5").Value = '  +66.79%  '
$ws.Range("D6").Value = '254.56'
$ws.Range("E6").Value = '  +3.23%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.3497'
$ws.Range("E8").Value = '  +20.88%  '
$ws.Range("D9").Value = '28.06'
$ws.Range("E9").Value = '  +26.23%  '
$ws.Range("D10").Value = '0.06996'
$ws.Range("E10").Value = '  +7.04%  '
$ws.Range("D11").Value = '0.8432'
$ws.Range("E11").Value = '  +8.61%  '
$ws.Range("D12").Value = '0.08194'
$ws.Range("E12").Value = '  +5.02%  '
$ws.Range("D13").Value = '1.991.37'
$ws.Range("E13").Value = '  +5.82%  '
$ws.Range("D14").Value = '100.24'
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '5.575'
$ws.Range("E15").Value = '  +6.11%  '
$ws.Range("D16").Value = '15.28'
$ws.Range("E16").Value = '  +15.66%  '
$ws.Range("D17").Value = '272.66'
$ws.Range("E17").Value = '  -4.52%  '
$ws.Range("D18").Value = '31.153.25'
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("D19").Value = '5.868'
$ws.Range("E19").Value = '  +9.40%  '
$ws.Range("D20").Value = '0.000007923'
$ws.Range("E20").Value = '  +5.25%  '
$ws.Range("D21").Value = '2.248.39'
$ws.Range("E21").Value = '  +5.80%  '
$ws.Range("D22").Value = '0.9993'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = '0.9994'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '7.039'
$ws.Range("E24").Value = '  +9.65%  '
$ws.Range("D25").Value = '9.977'
$ws.Range("E25").Value = '  +8.95%  '
$ws.Range("D26").Value = '0.1499'
$ws.Range("E26").Value = '  +54.49%  '
$ws.Range("D27").Value = '165.50'
$ws.Range("E27").Value = '  +1.42%  '
$ws.Range("D28").Value = '19.93'
$ws.Range("E28").Value = '  +4.39%  '
$ws.Range("D29").Value = '2.307'
$ws.Range("E29").Value = '  +20.38%  '
$ws.Range("D30").Value = '1.596'
$ws.Range("E30").Value = '  +6.00%  '
$ws.Range("D31").Value = '1.355'
$ws.Range("E31").Value = '  +1.70%  '
$ws.Range("D32").Value = '4.582'
$ws.Range("E32").Value = '  +7.44%  '
$ws.Range("D33").Value = '4.398'
$ws.Range("E33").Value = '  +4.85%  '
$ws.Range("D34").Value = '0.05220'
$ws.Range("E34").Value = '  +7.70%  '
$ws.Range("D35").Value = '1.225'
$ws.Range("E35").Value = '  +8.35%  '
$ws.Range("D36").Value = '0.7774'
$ws.Range("E36").Value = '  +11.42%  '
$ws.Range("D37").Value = '2.764'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02003'
$ws.Range("E38").Value = '  +4.29%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.887'
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '6.618'
$ws.Range("E40").Value = '  +5.16%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '79.15'
$ws.Range("E41").Value = '  +3.87%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.4656'
$ws.Range("E42").Value = '  +9.39%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '2.118'
$ws.Range("E43").Value = '  +6.48%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '105.22'
$ws.Range("E44").Value = '  +3.72%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.8492'
$ws.Range("E45").Value = '  +1.73%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '0.9996'
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.658'
$ws.Range("E47").Value = '  +8.73%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.861'
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '36.73'
$ws.Range("E49").Value = '  +4.21%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").Value = '0.4291'
$ws.Range("E50").Value = '  +8.48%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.520'
$ws.Range("E51").Value = '  +12.73%  '
